$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (existing rows 9..40 shift down to 10..41)
$ws.Rows("9:9").Insert()

# Populate the new row 9 with the weekly price-report entry
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44677
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107001
$ws.Range("J9").Value = "Caqui"
$ws.Range("K9").Value = "Fuyu"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 55
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "`$/bandeja 15 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 15
